$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 9 (shifts old rows 9,10,11 -> 10,11,12; the inserted
# row 9 initially has no borders/format)
$ws.Rows("9:9").Insert()

# --- Row 9 (new): replicate the data that used to live in row 8 ---
$ws.Range("A9").Value = "CodeGPTPy"
$ws.Range("B9").Value = 1024
$ws.Range("C9").Value = "SGD"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.675
$ws.Range("F9").Value = 0.71399999999999997
$ws.Range("G9").Formula = "=E9-F9"

# Match the border/number formatting used by the rest of the table
for ($c = 1; $c -le 7; $c++) {
  $ws.Cells.Item(9, $c).BorderAround(1, 2, 1, 0)
}
$ws.Range("E9:G9").NumberFormat = "0.000"

# --- Row 8: new data point (batch size 12288, Adam optimizer) ---
$ws.Range("B8").Value = 12288
$ws.Range("C8").Value = "Adam"
$ws.Range("E8").Value = 0.57999999999999996
$ws.Range("F8").Value = 0.67500000000000004
$ws.Range("G8").Formula = "=E8-F8"

# Extend the shared formula group that used to stop at G7 so it now
# also covers the (moved-down) G8
$ws.Range("G4").Formula = "=E4-F4"

# Reflect the last cell the author had selected when saving
$ws.Range("N40").Select()
